$d = $word.ActiveDocument

# Title paragraph: merge split runs "Questions:" " " "Introduction" " " "to" " " "hypothesis" " " "testing"
# into a single run "Questions: Introduction to hypothesis testing"
$d.Paragraphs(1).Range.Find.Execute(
    "Questions: Introduction to hypothesis testing", $false, $false, $false, $false, $false,
    $true, 1, $false, "Questions: Introduction to hypothesis testing", 2) | Out-Null

# Author paragraph: merge "Ellie" " " "Trace" into a single run "Ellie Trace"
$d.Paragraphs(2).Range.Find.Execute(
    "Ellie Trace", $false, $false, $false, $false, $false,
    $true, 1, $false, "Ellie Trace", 2) | Out-Null

# Abstract paragraph: merge all the split word/space runs into one run
$d.Paragraphs(4).Range.Find.Execute(
    "A selection of questions for the study guide on introduction to hypothesis testing.",
    $false, $false, $false, $false, $false,
    $true, 1, $false,
    "A selection of questions for the study guide on introduction to hypothesis testing.", 2) | Out-Null
